$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be auto-parsed as numbers
$textCells = @("D5", "D6", "D10", "D12", "D19", "D20", "D22", "D23", "D25", "D26", "D29", "D31", "D32", "D36", "D37", "D38", "D39", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values per the diff
$ws.Range("D2").Value = '62.357.93'
$ws.Range("E2").Value = '  -3.19%  '
$ws.Range("D3").Value = '3.369.20'
$ws.Range("E3").Value = '  -3.92%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '571.24'
$ws.Range("E5").Value = '  -3.40%  '
$ws.Range("D6").Value = '125.10'
$ws.Range("E6").Value = '  -7.08%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.368.20'
$ws.Range("E8").Value = '  -3.92%  '
$ws.Range("E9").Value = '  -2.64%  '
$ws.Range("D10").Value = '7.21'
$ws.Range("E10").Value = '  -5.60%  '
$ws.Range("E11").Value = '  -5.96%  '
$ws.Range("D12").Value = '0.374'
$ws.Range("E12").Value = '  -4.37%  '
$ws.Range("D13").Value = '3.949.83'
$ws.Range("E13").Value = '  -3.72%  '
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").Value = '3.372.95'
$ws.Range("E15").Value = '  -3.80%  '
$ws.Range("E16").Value = '  -6.99%  '
$ws.Range("D17").Value = '62.511.57'
$ws.Range("E17").Value = '  -2.94%  '
$ws.Range("E18").Value = '  -4.48%  '
$ws.Range("D19").Value = '9.08'
$ws.Range("E19").Value = '  -9.79%  '
$ws.Range("D20").Value = '5.55'
$ws.Range("E20").Value = '  -4.31%  '
$ws.Range("E21").Value = '  -4.09%  '
$ws.Range("D22").Value = '358.03'
$ws.Range("E22").Value = '  -8.66%  '
$ws.Range("D23").Value = '0.551'
$ws.Range("E23").Value = '  -5.70%  '
$ws.Range("D24").Value = '3.506.36'
$ws.Range("E24").Value = '  -3.82%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '70.94'
$ws.Range("E26").Value = '  -4.78%  '
$ws.Range("E27").Value = '  -10.70%  '
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").Value = '6.93'
$ws.Range("E29").Value = '  -7.24%  '
$ws.Range("E30").Value = '  -4.32%  '
$ws.Range("D31").Value = '7.80'
$ws.Range("E31").Value = '  -5.00%  '
$ws.Range("D32").Value = '2.10'
$ws.Range("E32").Value = '  -7.56%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '3.400.49'
$ws.Range("E34").Value = '  -3.79%  '
$ws.Range("E35").Value = '  -6.62%  '
$ws.Range("D36").Value = '22.49'
$ws.Range("E36").Value = '  -4.16%  '
$ws.Range("D37").Value = '5.32'
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("D38").Value = '167.48'
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("D39").Value = '6.60'
$ws.Range("E39").Value = '  -5.49%  '
$ws.Range("E40").Value = '  -5.91%  '
$ws.Range("D41").Value = '0.0751'
$ws.Range("E41").Value = '  -5.23%  '
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").Value = '41.65'
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("D44").Value = '0.761'
$ws.Range("E44").Value = '  -6.22%  '
$ws.Range("D45").Value = '4.19'
$ws.Range("E45").Value = '  -6.19%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").Value = '1.52'
$ws.Range("E46").Value = '  -8.53%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = '1.08'
$ws.Range("E47").Value = '  -8.89%  '
$ws.Range("D48").Value = '22.26'
$ws.Range("E48").Value = '  -10.76%  '
$ws.Range("D49").Value = '6.55'
$ws.Range("E49").Value = '  -4.08%  '
$ws.Range("D50").Value = '2.216.75'
$ws.Range("E50").Value = '  -8.24%  '
$ws.Range("D51").Value = '0.834'
$ws.Range("E51").Value = '  -10.37%  '

Write-Host "Applied cryptos update"
